# Generate Report for Handoff
#
# Adds a new localization entry for file 9faa18dc-475c-4d23-9dd8-37a8a75cf24c.md
# as row 3 on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the
# existing row 2 pattern that was recorded for
# 03362531-224f-4491-b911-5d6025e47ac8.md.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0333a6aeee5c19a42226f4f4c27ed044e064d20/e2e/"
$fileGuid = "9faa18dc-475c-4d23-9dd8-37a8a75cf24c"
$mdName = "$fileGuid.md"
$mdPath = "e2e\$mdName"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.ListObjects.Item(1).ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "'$mdName"
$wsOverview.Range("C3").Value = "'.md"
$wsOverview.Range("D3").Value = "'"
$wsOverview.Range("E3").Value = "'Ready for handoff"
$wsOverview.Range("F3").Value = "'Ready for handoff"
$wsOverview.Range("G3").Value = "'2016-09-01 18:46:54"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$baseUrl$mdName", "", "", $mdPath)
$wsOverview.Range("B3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.ListObjects.Item(1).ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = "'.md"
$wsZhCn.Range("C3").Value = "'Ready for handoff"
$wsZhCn.Range("D3").Value = "'e2e"
$wsZhCn.Range("E3").Value = "'ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "'$fileGuid.f86d78c7062d5debc1f94262a7998dc2c25b3bcf.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "'2016-09-01 18:46:50"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("I3").Value = "'"
$wsZhCn.Range("J3").Value = "'"
$wsZhCn.Range("K3").Value = "'0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "$baseUrl$mdName", "", "", $mdName)
$wsZhCn.Range("A3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.ListObjects.Item(1).ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = "'.md"
$wsDeDe.Range("C3").Value = "'Ready for handoff"
$wsDeDe.Range("D3").Value = "'e2e"
$wsDeDe.Range("E3").Value = "'ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "'$fileGuid.f86d78c7062d5debc1f94262a7998dc2c25b3bcf.de-de.xlf"
$wsDeDe.Range("H3").Value = "'2016-09-01 18:46:54"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("I3").Value = "'"
$wsDeDe.Range("J3").Value = "'"
$wsDeDe.Range("K3").Value = "'0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "$baseUrl$mdName", "", "", $mdName)
$wsDeDe.Range("A3").Style = "HyperLink"
